# Rename the "T_Shirts" worksheet to "T-Shirts"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("T_Shirts")
$ws.Name = "T-Shirts"
